# Discharge_July26.xlsx - "lots of discharge data"
# Adds a new "New depth" measurement block (rows 34-47) to the "stn3" sheet,
# mirroring the structure of the existing X/V/D/segment/Q/Qtotal block above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stn3")

# --- Row 34: bold section header ---------------------------------------
$ws.Range("A34").Value = "New depth"
$ws.Range("A34").Font.Bold = $true

# --- Row 35: column headers (reuse the same labels as row 18) ----------
$ws.Range("A35").Value = "X"
$ws.Range("B35").Value = "V"
$ws.Range("C35").Value = "D"
$ws.Range("D35").Value = "segment"
$ws.Range("E35").Value = "Q"
$ws.Range("F35").Value = "Qtotal"

# --- Column A: station positions (plain values) -------------------------
$ws.Range("A36").Value = 0.75
$ws.Range("A37").Value = 0.8
$ws.Range("A38").Value = 0.85
$ws.Range("A39").Value = 0.9
$ws.Range("A40").Value = 0.95
$ws.Range("A41").Value = 1
$ws.Range("A42").Value = 1.05
$ws.Range("A43").Value = 1.1
$ws.Range("A44").Value = 1.15
$ws.Range("A45").Value = 1.2
$ws.Range("A46").Value = 1.25
$ws.Range("A47").Value = 1.3

# --- Column B: velocities -- computed the same way as the B20:B30 block
# (0.0572 * the corresponding row from the first B4:B14 block), then
# pasted as plain values (matches the workbook: B37:B47 hold <v> only).
$bMap = @{
  37 = 4
  38 = 5
  39 = 6
  40 = 7
  41 = 8
  42 = 9
  43 = 10
  44 = 11
  45 = 12
  46 = 13
  47 = 14
}
$ws.Range("B36").Value = 0
foreach ($row in 37..47) {
  $src = $bMap[$row]
  $cell = $ws.Range("B$row")
  $cell.Formula = "=0.0572*B$src"
  $v = $cell.Value2
  $cell.Value = $v
}

# --- Column C: depth in cm, converted from the D column of the block above
$ws.Range("C36").Formula = "=C19*2.54"
$ws.Range("C37:C47").Formula = "=C20*2.54"

# --- Column D: midpoint segment positions --------------------------------
$ws.Range("D36").Formula = "=A36"
$ws.Range("D37").Formula = "=(A37+(A38-A37)/2)"
$ws.Range("D38:D39").Formula = "=(A38+(A39-A38)/2)"
$ws.Range("D40").Formula = "=(A40+(A41-A40)/2)"
$ws.Range("D41:D47").Formula = "=(A41+(A42-A41)/2)"

# --- Column E: segment discharge -----------------------------------------
$ws.Range("E37").Formula = "=(D37-D36)*(B37)*C37"
$ws.Range("E38").Formula = "=(D38-D37)*(B38)*C38"
$ws.Range("E39:E47").Formula = "=(D39-D38)*(B39)*C39"

# --- Column F: total discharge for the block -----------------------------
$ws.Range("F36").Formula = "=SUM(E36:E54)"

# --- View state: "stn3" becomes the active sheet/tab, with F36 selected --
$ws.Activate()
$ws.Range("F36").Select()
